# Fix duplicate player names by normalizing diacritics (Jokic/Jokić -> Jokic,
# Doncic/Dončić -> Doncic). Merging the duplicate rows shifts several
# leaderboard rows so names/scores are realigned below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-12: "Nikola Jokic" (with diacritic duplicate removed) moves to the
# top of its score tier (11 games); the rest of the tier shifts down one row.
$ws.Range("A8").Value = 'Nikola'
$ws.Range("B8").Value = 'Jokic'
$ws.Range("C8").Value = 11
$ws.Range("A9").Value = 'Ben'
$ws.Range("B9").Value = 'Wallace'
$ws.Range("C9").Value = 10
$ws.Range("A10").Value = 'Marcus'
$ws.Range("B10").Value = 'Camby'
$ws.Range("A11").Value = 'Lafayette'
$ws.Range("B11").Value = 'Lever'
$ws.Range("A12").Value = 'Anthony'
$ws.Range("B12").Value = 'Davis'

# Rows 29-42: "Luka Doncic" merges up into the 4-game tier, bumping George
# McGinnis/Clyde Drexler's tier accordingly, and shifting neighbors by one row.
$ws.Range("A29").Value = 'Larry'
$ws.Range("B29").Value = 'Steele'
$ws.Range("A34").Value = 'Luka'
$ws.Range("B34").Value = 'Doncic'
$ws.Range("A36").Value = 'Mark'
$ws.Range("B36").Value = 'Eaton'
$ws.Range("A37").Value = 'Gerald'
$ws.Range("B37").Value = 'Wallace'
$ws.Range("A38").Value = 'Giannis'
$ws.Range("B38").Value = 'Antetokounmpo'
$ws.Range("A39").Value = 'Bob'
$ws.Range("B39").Value = 'McAdoo'
$ws.Range("A40").Value = 'Clyde'
$ws.Range("B40").Value = 'Drexler'
$ws.Range("C40").Value = 4
$ws.Range("A42").Value = 'George'
$ws.Range("B42").Value = 'McGinnis'

# Rows 89-90: swap order of Damian Lillard / Charles Jones.
$ws.Range("A89").Value = 'Charles'
$ws.Range("B89").Value = 'Jones'
$ws.Range("A90").Value = 'Damian'
$ws.Range("B90").Value = 'Lillard'

# Rows 101-102: swap order of Quinn Buckner / Paul Pressey.
$ws.Range("A101").Value = 'Paul'
$ws.Range("B101").Value = 'Pressey'
$ws.Range("A102").Value = 'Quinn'
$ws.Range("B102").Value = 'Buckner'

# Rows 107-117: the diacritic duplicate "Nikola Jokić" (1 game) is removed
# from row 107 (merged into row 8 above), so rows 108-117 shift up by one,
# and a new unique player ("Moussa Diabate") is appended at row 117 to keep
# the row count/layout intact.
$ws.Range("A107").Value = 'Brook'
$ws.Range("B107").Value = 'Lopez'
$ws.Range("A112").Value = 'Nicolas'
$ws.Range("B112").Value = 'Batum'
$ws.Range("A113").Value = 'Nick'
$ws.Range("B113").Value = 'Anderson'
$ws.Range("A114").Value = 'Cade'
$ws.Range("B114").Value = 'Cunningham'
$ws.Range("A115").Value = 'Charles'
$ws.Range("B115").Value = 'Oakley'
$ws.Range("B116").Value = 'Smith'
$ws.Range("A117").Value = 'Moussa'
$ws.Range("B117").Value = 'Diabate'
